# Refresh the cryptos sheet with newly scraped Price / Volume(1h)
# figures (GitHub Actions data-refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.729.35"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").Value = "3.786.19"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("D5").Value = "597.11"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").Value = "169.57"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").Value = "3.785.87"
$ws.Range("E7").Value = "  -1.79%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("E10").Value = "  -0.77%  "

$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").Value = "0.0000279"
$ws.Range("E13").Value = "  +4.59%  "

$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").Value = "4.421.33"
$ws.Range("E15").Value = "  -1.82%  "

$ws.Range("D16").Value = "3.789.64"
$ws.Range("E16").Value = "  -1.75%  "

$ws.Range("D17").Value = "18.57"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").Value = "67.747.71"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("E20").Value = "  +0.73%  "

$ws.Range("E21").Value = "  -6.14%  "

$ws.Range("D22").Value = "469.06"
$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").Value = "0.719"
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("E24").Value = "  -7.89%  "

$ws.Range("D25").NumberFormat = "@"  # keep trailing zero as text
$ws.Range("D25").Value = "83.90"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("E26").Value = "  -1.54%  "

$ws.Range("E27").Value = "  +0.42%  "

$ws.Range("D28").Value = "10.29"
$ws.Range("E28").Value = "  +0.74%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "2.92"
$ws.Range("E30").Value = "  -1.73%  "

$ws.Range("D31").Value = "3.935.12"
$ws.Range("E31").Value = "  -1.78%  "

$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("D33").Value = "30.58"
$ws.Range("E33").Value = "  -2.95%  "

$ws.Range("E34").Value = "  -3.59%  "

$ws.Range("D35").Value = "9.16"
$ws.Range("E35").Value = "  -2.13%  "

$ws.Range("D36").Value = "3.746.65"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("E37").Value = "  +1.79%  "

$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("E39").Value = "  -1.29%  "

$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("E41").Value = "  -2.23%  "

$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.15%  "

$ws.Range("E43").Value = "  -1.11%  "

$ws.Range("D45").Value = "8.69"
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("E46").Value = "  -1.99%  "

$ws.Range("E47").Value = "  -2.75%  "

$ws.Range("D48").NumberFormat = "@"  # keep trailing zero as text
$ws.Range("D48").Value = "395.60"
$ws.Range("E48").Value = "  -5.11%  "

$ws.Range("E49").Value = "  -8.09%  "

$ws.Range("E50").Value = "  -1.02%  "

$ws.Range("D51").Value = "39.27"
$ws.Range("E51").Value = "  +3.40%  "
